$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Narucitelj" block rework.
#    Before: an empty paragraph (carrying a stray red/FF0000 rPr on its
#    mark) immediately followed by a single paragraph that reads
#    "Naručitelj: ${IME} ${PREZIME}, ${ADRESA} OIB: ${OIB}".
#    After: the empty paragraph is gone, "Naručitelj: " becomes its own
#    paragraph, and the placeholder list collapses to "${NARUCITELJI}"
#    living alone (no longer justified) in the paragraph that follows.
# ---------------------------------------------------------------------

# 1a) Drop the stray empty paragraph right above "Naručitelj: ".
$find1 = $d.Content
$find1.Find.Execute('Naručitelj: ', $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$naruciteljIndex = $find1.Paragraphs.Item(1).Index
$prevPara = $d.Paragraphs.Item($naruciteljIndex - 1)
$prevPara.Range.Delete()

# 1b) Split "Naručitelj: " into its own paragraph.
$c = $d.Content
$c.Find.Execute('Naručitelj: ', $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$c.InsertParagraphAfter()

# 1c) Collapse "${IME} ${PREZIME}, ${ADRESA} OIB: ${OIB}" to "${NARUCITELJI}".
$r1 = $d.Content
$r1.Find.Execute('${IME}', $false, $false, $false, $false, $false, $true, 1, $false, '${', 2) | Out-Null

$r2 = $d.Content
$r2.Find.Execute(' ${PREZIME}', $false, $false, $false, $false, $false, $true, 1, $false, 'NARUCITELJI', 2) | Out-Null

$r3 = $d.Content
$r3.Find.Execute(', ${ADRESA} OIB: ', $false, $false, $false, $false, $false, $true, 1, $false, '', 2) | Out-Null

$r4 = $d.Content
$r4.Find.Execute('${OIB}', $false, $false, $false, $false, $false, $true, 1, $false, '}', 2) | Out-Null

# 1d) The "${NARUCITELJI}" paragraph is no longer justified ("both").
$r5 = $d.Content
$r5.Find.Execute('${NARUCITELJI}', $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$narucIndex = $r5.Paragraphs.Item(1).Index
$d.Paragraphs.Item($narucIndex).Alignment = 0

# ---------------------------------------------------------------------
# 2) Remove one of the two blank (sz 24) paragraphs that sit right
#    after the tab-only paragraph below "Oznaka elaborata: ...".
# ---------------------------------------------------------------------
$tabFind = $d.Content
$tabFind.Find.Execute("Oznaka elaborata: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$oznakaIndex = $tabFind.Paragraphs.Item(1).Index
$blankIndex = $oznakaIndex + 2
$d.Paragraphs.Item($blankIndex).Range.Delete()
